$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 3) mirroring row 2's values, but with a new
# "phoneNumber"-style entry in column A (plain text, no hyperlink this time)
# and no value in column H.
$ws.Range("A3").Value = "Mani12**("
$ws.Range("B3").Value = "Raju"
$ws.Range("C3").Value = "Devathi"
$ws.Range("D3").Value = """20-08-1993"""
$ws.Range("E3").Value = "Male"
$ws.Range("F3").Value = "Bangalore"
$ws.Range("G3").Value = "I have an active credit card"

# Match D2's date-style formatting on D3 by copying just its format.
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the selected cell to reflect the new extent of used data.
$ws.Range("H3").Select()
